# Update the cryptocurrency ranking table (columns Price/Volume(1h), and the
# Coin/Link pair for rows 35-36 which swapped rank position) to the latest
# scrape snapshot.
#
# Column D ("Price") holds values that LOOK numeric (e.g. "576.98") but must
# stay plain text, matching the source sheet's inlineStr cells -- otherwise
# Excel auto-converts the assignment to a float and the "." thousand
# separators used elsewhere in the column (e.g. "63.019.66") would become
# inconsistent with this column's type. Setting NumberFormat to "@" (Text)
# before writing to cells whose new value parses as a number keeps them text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '63.019.66'
$ws.Cells.Item(2, 5).Value = '  -1.23%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '2.549.88'
$ws.Cells.Item(3, 5).Value = '  -0.53%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '576.98'
$ws.Cells.Item(5, 5).Value = '  +0.22%  '
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.79'
$ws.Cells.Item(6, 5).Value = '  -2.99%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.05%  '
# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.582'
$ws.Cells.Item(8, 5).Value = '  -0.98%  '
# Row 9
$ws.Cells.Item(9, 5).Value = '  -1.58%  '
# Row 10
$ws.Cells.Item(10, 5).Value = '  -4.98%  '
# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.64%  '
# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.51%  '
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '27.14'
$ws.Cells.Item(13, 5).Value = '  -4.12%  '
# Row 14
$ws.Cells.Item(14, 4).Value = '3.006.23'
$ws.Cells.Item(14, 5).Value = '  -0.48%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '62.973.15'
$ws.Cells.Item(15, 5).Value = '  -1.08%  '
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000142'
$ws.Cells.Item(16, 5).Value = '  -1.41%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '2.552.58'
$ws.Cells.Item(17, 5).Value = '  -0.23%  '
# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '11.33'
$ws.Cells.Item(18, 5).Value = '  -2.63%  '
# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '335.38'
$ws.Cells.Item(19, 5).Value = '  -2.32%  '
# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.24%  '
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.74'
$ws.Cells.Item(21, 5).Value = '  -2.60%  '
# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.13%  '
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '65.36'
$ws.Cells.Item(23, 5).Value = '  -1.29%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.79%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.85%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.06%  '
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.48'
$ws.Cells.Item(27, 5).Value = '  +3.31%  '
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.32'
$ws.Cells.Item(28, 5).Value = '  -1.83%  '
# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.34'
$ws.Cells.Item(29, 5).Value = '  +2.52%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.39%  '
# Row 31
$ws.Cells.Item(31, 4).Value = '0.0₃0811'
$ws.Cells.Item(31, 5).Value = '  -4.30%  '
# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '178.05'
$ws.Cells.Item(32, 5).Value = '  +0.40%  '
# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.54'
$ws.Cells.Item(33, 5).Value = '  -3.96%  '
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '406.17'
$ws.Cells.Item(34, 5).Value = '  -4.19%  '
# Row 35
$ws.Cells.Item(35, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.399'
$ws.Cells.Item(35, 5).Value = '  -2.22%  '
# Row 36
$ws.Cells.Item(36, 2).Value = 'EthereumClassic'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '19.09'
$ws.Cells.Item(36, 5).Value = '  -0.50%  '
# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '4.33'
$ws.Cells.Item(38, 5).Value = '  -3.72%  '
# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.73'
$ws.Cells.Item(39, 5).Value = '  -2.24%  '
# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.10%  '
# Row 41
$ws.Cells.Item(41, 5).Value = '  -2.35%  '
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '151.26'
$ws.Cells.Item(42, 5).Value = '  -3.71%  '
# Row 43
$ws.Cells.Item(43, 5).Value = '  -1.74%  '
# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '20.84'
$ws.Cells.Item(44, 5).Value = '  -1.56%  '
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0538'
$ws.Cells.Item(45, 5).Value = '  +0.50%  '
# Row 46
$ws.Cells.Item(46, 5).Value = '  -1.89%  '
# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.51%  '
# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0237'
$ws.Cells.Item(48, 5).Value = '  +1.40%  '
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '18.24'
$ws.Cells.Item(49, 5).Value = '  -3.95%  '
# Row 50
$ws.Cells.Item(50, 5).Value = '  -10.00%  '
# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.35%  '
